$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the data row (row 2): Name, file_date, is_random
# A2: Name value "جيانا" -> "ديما"
$ws.Range("A2").Value = "ديما"
# D2: is_random "0" -> "1"  (set before C2 so shared-string order matches)
$ws.Range("D2").Value = "1"
# C2: file_date "02-03" -> "02-02"
$ws.Range("C2").Value = "02-02"

# Move the active selection to B4
$ws.Range("B4").Select()

$wb.Save()
